$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Pecam1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 52.61615766666667
$ws.Cells.Item(2, 8).Value = 157.848473
$ws.Cells.Item(2, 9).Value = 0.7671520491359202
$ws.Cells.Item(2, 10).Value = 0.7671520491359202
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 343.9479473333333
$ws.Cells.Item(2, 14).Value = 1031.843842
$ws.Cells.Item(2, 15).Value = 0.9666099193889262
$ws.Cells.Item(2, 16).Value = 0.966609919388926
$ws.Cells.Item(2, 17).Value = 18097.21942601703
$ws.Cells.Item(2, 18).Value = 162874.9748341533
$ws.Cells.Item(2, 19).Value = 0.7415367803743214
$ws.Cells.Item(2, 20).Value = 0.7415367803743212

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Pecam1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 52.61615766666667
$ws.Cells.Item(3, 8).Value = 157.848473
$ws.Cells.Item(3, 9).Value = 0.7671520491359202
$ws.Cells.Item(3, 10).Value = 0.7671520491359202
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.970184
$ws.Cells.Item(3, 14).Value = 17.910552
$ws.Cells.Item(3, 15).Value = 0.01677823379880302
$ws.Cells.Item(3, 16).Value = 0.01677823379880302
$ws.Cells.Item(3, 17).Value = 314.1281426430107
$ws.Cells.Item(3, 18).Value = 2827.153283787096
$ws.Cells.Item(3, 19).Value = 0.01287145643963329
$ws.Cells.Item(3, 20).Value = 0.01287145643963329

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Pecam1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 52.61615766666667
$ws.Cells.Item(4, 8).Value = 157.848473
$ws.Cells.Item(4, 9).Value = 0.7671520491359202
$ws.Cells.Item(4, 10).Value = 0.7671520491359202
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.671367333333334
$ws.Cells.Item(4, 14).Value = 8.014102000000001
$ws.Cells.Item(4, 15).Value = 0.007507444608265281
$ws.Cells.Item(4, 16).Value = 0.00750744460826528
$ws.Cells.Item(4, 17).Value = 140.5570847962496
$ws.Cells.Item(4, 18).Value = 1265.013763166246
$ws.Cells.Item(4, 19).Value = 0.005759351515005126
$ws.Cells.Item(4, 20).Value = 0.005759351515005126

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Pecam1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 52.61615766666667
$ws.Cells.Item(5, 8).Value = 157.848473
$ws.Cells.Item(5, 9).Value = 0.7671520491359202
$ws.Cells.Item(5, 10).Value = 0.7671520491359202
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.239611333333334
$ws.Cells.Item(5, 14).Value = 9.718834000000001
$ws.Cells.Item(5, 15).Value = 0.009104402204005551
$ws.Cells.Item(5, 16).Value = 0.00910440220400555
$ws.Cells.Item(5, 17).Value = 170.455900693387
$ws.Cells.Item(5, 18).Value = 1534.103106240482
$ws.Cells.Item(5, 19).Value = 0.006984460806960447
$ws.Cells.Item(5, 20).Value = 0.006984460806960445

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Pecam1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.395935333333332
$ws.Cells.Item(6, 8).Value = 28.187806
$ws.Cells.Item(6, 9).Value = 0.1369942497546098
$ws.Cells.Item(6, 10).Value = 0.1369942497546098
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 343.9479473333333
$ws.Cells.Item(6, 14).Value = 1031.843842
$ws.Cells.Item(6, 15).Value = 0.9666099193889262
$ws.Cells.Item(6, 16).Value = 0.966609919388926
$ws.Cells.Item(6, 17).Value = 3231.712671176739
$ws.Cells.Item(6, 18).Value = 29085.41404059065
$ws.Cells.Item(6, 19).Value = 0.1324200007120498
$ws.Cells.Item(6, 20).Value = 0.1324200007120498

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Pecam1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.395935333333332
$ws.Cells.Item(7, 8).Value = 28.187806
$ws.Cells.Item(7, 9).Value = 0.1369942497546098
$ws.Cells.Item(7, 10).Value = 0.1369942497546098
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.970184
$ws.Cells.Item(7, 14).Value = 17.910552
$ws.Cells.Item(7, 15).Value = 0.01677823379880302
$ws.Cells.Item(7, 16).Value = 0.01677823379880302
$ws.Cells.Item(7, 17).Value = 56.09546279210132
$ws.Cells.Item(7, 18).Value = 504.859165128912
$ws.Cells.Item(7, 19).Value = 0.002298521551474457
$ws.Cells.Item(7, 20).Value = 0.002298521551474457

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Pecam1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.395935333333332
$ws.Cells.Item(8, 8).Value = 28.187806
$ws.Cells.Item(8, 9).Value = 0.1369942497546098
$ws.Cells.Item(8, 10).Value = 0.1369942497546098
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.671367333333334
$ws.Cells.Item(8, 14).Value = 8.014102000000001
$ws.Cells.Item(8, 15).Value = 0.007507444608265281
$ws.Cells.Item(8, 16).Value = 0.00750744460826528
$ws.Cells.Item(8, 17).Value = 25.09999471557911
$ws.Cells.Item(8, 18).Value = 225.899952440212
$ws.Cells.Item(8, 19).Value = 0.001028476741683593
$ws.Cells.Item(8, 20).Value = 0.001028476741683593

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Pecam1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.395935333333332
$ws.Cells.Item(9, 8).Value = 28.187806
$ws.Cells.Item(9, 9).Value = 0.1369942497546098
$ws.Cells.Item(9, 10).Value = 0.1369942497546098
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.239611333333334
$ws.Cells.Item(9, 14).Value = 9.718834000000001
$ws.Cells.Item(9, 15).Value = 0.009104402204005551
$ws.Cells.Item(9, 16).Value = 0.00910440220400555
$ws.Cells.Item(9, 17).Value = 30.43917859313378
$ws.Cells.Item(9, 18).Value = 273.952607338204
$ws.Cells.Item(9, 19).Value = 0.001247250749401956
$ws.Cells.Item(9, 20).Value = 0.001247250749401956

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Pecam1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.08161033333333334
$ws.Cells.Item(10, 8).Value = 0.244831
$ws.Cells.Item(10, 9).Value = 0.001189891797952309
$ws.Cells.Item(10, 10).Value = 0.001189891797952309
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 343.9479473333333
$ws.Cells.Item(10, 14).Value = 1031.843842
$ws.Cells.Item(10, 15).Value = 0.9666099193889262
$ws.Cells.Item(10, 16).Value = 0.966609919388926
$ws.Cells.Item(10, 17).Value = 28.06970663118911
$ws.Cells.Item(10, 18).Value = 252.627359680702
$ws.Cells.Item(10, 19).Value = 0.001150161214900226
$ws.Cells.Item(10, 20).Value = 0.001150161214900225

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Efnb2"
$ws.Cells.Item(11, 3).Value = "Pecam1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.08161033333333334
$ws.Cells.Item(11, 8).Value = 0.244831
$ws.Cells.Item(11, 9).Value = 0.001189891797952309
$ws.Cells.Item(11, 10).Value = 0.001189891797952309
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5.970184
$ws.Cells.Item(11, 14).Value = 17.910552
$ws.Cells.Item(11, 15).Value = 0.01677823379880302
$ws.Cells.Item(11, 16).Value = 0.01677823379880302
$ws.Cells.Item(11, 17).Value = 0.4872287063013334
$ws.Cells.Item(11, 18).Value = 4.385058356712
$ws.Cells.Item(11, 19).Value = 0.00001996428278132193
$ws.Cells.Item(11, 20).Value = 0.00001996428278132193

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Efnb2"
$ws.Cells.Item(12, 3).Value = "Pecam1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.08161033333333334
$ws.Cells.Item(12, 8).Value = 0.244831
$ws.Cells.Item(12, 9).Value = 0.001189891797952309
$ws.Cells.Item(12, 10).Value = 0.001189891797952309
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.671367333333334
$ws.Cells.Item(12, 14).Value = 8.014102000000001
$ws.Cells.Item(12, 15).Value = 0.007507444608265281
$ws.Cells.Item(12, 16).Value = 0.00750744460826528
$ws.Cells.Item(12, 17).Value = 0.2180111785291112
$ws.Cells.Item(12, 18).Value = 1.962100606762
$ws.Cells.Item(12, 19).Value = 0.000008933046762956142
$ws.Cells.Item(12, 20).Value = 0.000008933046762956142

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Efnb2"
$ws.Cells.Item(13, 3).Value = "Pecam1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.08161033333333334
$ws.Cells.Item(13, 8).Value = 0.244831
$ws.Cells.Item(13, 9).Value = 0.001189891797952309
$ws.Cells.Item(13, 10).Value = 0.001189891797952309
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.239611333333334
$ws.Cells.Item(13, 14).Value = 9.718834000000001
$ws.Cells.Item(13, 15).Value = 0.009104402204005551
$ws.Cells.Item(13, 16).Value = 0.00910440220400555
$ws.Cells.Item(13, 17).Value = 0.2643857607837778
$ws.Cells.Item(13, 18).Value = 2.379471847054
$ws.Cells.Item(13, 19).Value = 0.00001083325350780513
$ws.Cells.Item(13, 20).Value = 0.00001083325350780513

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Efnb2"
$ws.Cells.Item(14, 3).Value = "Pecam1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 6.492645
$ws.Cells.Item(14, 8).Value = 19.477935
$ws.Cells.Item(14, 9).Value = 0.09466380931151776
$ws.Cells.Item(14, 10).Value = 0.09466380931151776
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 343.9479473333333
$ws.Cells.Item(14, 14).Value = 1031.843842
$ws.Cells.Item(14, 15).Value = 0.9666099193889262
$ws.Cells.Item(14, 16).Value = 0.966609919388926
$ws.Cells.Item(14, 17).Value = 2233.13192051403
$ws.Cells.Item(14, 18).Value = 20098.18728462627
$ws.Cells.Item(14, 19).Value = 0.09150297708765485
$ws.Cells.Item(14, 20).Value = 0.09150297708765484

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Efnb2"
$ws.Cells.Item(15, 3).Value = "Pecam1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 6.492645
$ws.Cells.Item(15, 8).Value = 19.477935
$ws.Cells.Item(15, 9).Value = 0.09466380931151776
$ws.Cells.Item(15, 10).Value = 0.09466380931151776
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 5.970184
$ws.Cells.Item(15, 14).Value = 17.910552
$ws.Cells.Item(15, 15).Value = 0.01677823379880302
$ws.Cells.Item(15, 16).Value = 0.01677823379880302
$ws.Cells.Item(15, 17).Value = 38.76228529668
$ws.Cells.Item(15, 18).Value = 348.86056767012
$ws.Cells.Item(15, 19).Value = 0.001588291524913952
$ws.Cells.Item(15, 20).Value = 0.001588291524913952

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Efnb2"
$ws.Cells.Item(16, 3).Value = "Pecam1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 6.492645
$ws.Cells.Item(16, 8).Value = 19.477935
$ws.Cells.Item(16, 9).Value = 0.09466380931151776
$ws.Cells.Item(16, 10).Value = 0.09466380931151776
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 2.671367333333334
$ws.Cells.Item(16, 14).Value = 8.014102000000001
$ws.Cells.Item(16, 15).Value = 0.007507444608265281
$ws.Cells.Item(16, 16).Value = 0.00750744460826528
$ws.Cells.Item(16, 17).Value = 17.34423975993001
$ws.Cells.Item(16, 18).Value = 156.09815783937
$ws.Cells.Item(16, 19).Value = 0.0007106833048136067
$ws.Cells.Item(16, 20).Value = 0.0007106833048136067

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Efnb2"
$ws.Cells.Item(17, 3).Value = "Pecam1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 6.492645
$ws.Cells.Item(17, 8).Value = 19.477935
$ws.Cells.Item(17, 9).Value = 0.09466380931151776
$ws.Cells.Item(17, 10).Value = 0.09466380931151776
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.239611333333334
$ws.Cells.Item(17, 14).Value = 9.718834000000001
$ws.Cells.Item(17, 15).Value = 0.009104402204005551
$ws.Cells.Item(17, 16).Value = 0.00910440220400555
$ws.Cells.Item(17, 17).Value = 21.03364632531001
$ws.Cells.Item(17, 18).Value = 189.30281692779
$ws.Cells.Item(17, 19).Value = 0.0008618573941353436
$ws.Cells.Item(17, 20).Value = 0.0008618573941353433

